# daily auto push: 2026-01-07 02:26 UTC
# Insert 4 new hourly readings for 2026/01/07 (Wednesday) right after the
# existing 2026/01/07 03:00 row (row 580), shifting all subsequent rows
# down by four positions. This mirrors the upstream diff, which inserts
# rows for hours 7, 8, 9 and 10 on 2026/01/07 and pushes the remaining
# (2026/12/29 .. 2027/01/05) data down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 581 (pushes old rows 581:622 -> 585:626)
$ws.Rows("581:584").Insert()

$newRows = @(
    @{Row = 581; Date = "2026/01/07"; Day = "水"; Hour = 7},
    @{Row = 582; Date = "2026/01/07"; Day = "水"; Hour = 8},
    @{Row = 583; Date = "2026/01/07"; Day = "水"; Hour = 9},
    @{Row = 584; Date = "2026/01/07"; Day = "水"; Hour = 10}
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Force column A to be stored as plain text instead of Excel
    # auto-detecting/parsing the "YYYY/MM/DD" string as a date value.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $item.Date
    # Drop the temporary text-format styling so the cell keeps the same
    # (default/no explicit style) appearance as the rest of the data rows.
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 2).Value = $item.Day
    $ws.Cells.Item($r, 3).Value = $item.Hour
    $ws.Cells.Item($r, 4).Value = 201
}
